$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing URL for the existing Dalhousie University row
$ws.Range("C14").Value = "http://libraries.dal.ca/"

# Add new row 15: Dartmouth College
$ws.Range("A15").Value = "Dartmouth College"
$ws.Range("B15").Value = "Dartmouth Biomedical Libraries"
$ws.Range("C15").Value = "http://www.dartmouth.edu/~library/biomed/?mswitch-redir=classic"

# Add new row 16: Drexel University Libraries
$ws.Range("A16").Value = "Drexel University Libraries"
$ws.Range("B16").Value = "Library Academic Partnerships"

# Wrap text on the institution-name cells of the two new rows
$ws.Range("A15:A16").WrapText = $true

# Select the final cell, mimicking where the author left off editing
$ws.Range("C16").Select()
